$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, `
                         $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null
}

# Application/contract number "2404367" -> "1122334" (appears twice in the document)
Replace-Text "2404367" "1122334"

# Contract date "06.05.2024" -> "07.03.2024" (appears twice in the document)
Replace-Text "06.05.2024" "07.03.2024"

# Vessel name "СУВОРОВЕЦ" -> "СИНЕГОРСК"
Replace-Text "СУВОРОВЕЦ" "СИНЕГОРСК"

# Registry number "802465" -> "021026"
Replace-Text "802465" "021026"

# Survey description
Replace-Text "Первоначальное освидетельствование маломерного судна" "Первоначальное освидетельствование"

# Certificate description
Replace-Text "Удостоверение ф. 6.3.80 № 24.43.02.00456.121 от 07.05.2024" "Свидетельство ф. 8.5.3 № 24.42.02.00123.121 от 04.05.2024"

# Cost amount
Replace-Text "1 000,00 p. (одна тысяча рублей 00 копеек)" "100 000,00 p. (сто тысяч рублей 00 копеек)"
